# Adding svat_ok_activity_start to i18n
# Insert a new row in the "i18n" table just above the "tax_registration_number"
# row (row 542) holding the new translation key + its Portuguese text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 542 - this shifts every row below it down by
# one and naturally inherits the formatting (styles) of the row above, just
# like Excel does when a user inserts a row by hand.
$ws.Rows(542).Insert()

# Row height used by this table for the wrapped-text translation rows.
$ws.Rows(542).RowHeight = 34

# Fill in the new key/value pair.
$ws.Range("A542").Value = "svat_ok_activity_start"
$ws.Range("B542").Value = "Não foi efetuado este teste uma vez que a empresa teve o início de atividade em {0}."

# Grow the "i18n" table definition so it covers the freshly inserted row.
$lo = $ws.ListObjects.Item("i18n")
$lo.Resize($ws.Range("A1:G568"))

# Match the author's final on-screen selection/scroll state.
$win = $excel.ActiveWindow
$win.ScrollRow = 528
$win.ScrollColumn = 1
$ws.Range("A543").Select()
